$wb = $excel.ActiveWorkbook

# New company "bgwrer" (active = TRUE) added to the Companies lookup sheet.
$companies = $wb.Worksheets.Item("Companies")
$companies.Range("A3").Value = "bgwrer"
$companies.Range("B3").Value = "'TRUE"

# The existing "BC" location now also maps to the new "bgwrer" company.
$locations = $wb.Worksheets.Item("Locations")
$locations.Range("A3").Value = "BC"
$locations.Range("B3").Value = "bgwrer"
